$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This text is shared across the Overview sheet (B/C columns) and the
#    per-locale "Status" column (C) on zh-cn / de-de, so every cell that
#    currently shows it needs to be refreshed to keep them in sync.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. Handback report columns: fill in "Latest Target File" (F) and
#    "Latest Handback File" (G) for both locale sheets, rows 2 & 3, with
#    hyperlinks mirroring the existing Source File (A/B) and Handoff File
#    (D) hyperlinks.
# ---------------------------------------------------------------------------

# zh-cn -----------------------------------------------------------------
$zhcn.Range("F2").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b3af6e214d11ac0ecaa6b80385723228932fa6be/e2e/430ce049-5402-4497-a1f8-90bcfc8e1e5c.md", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md") | Out-Null

$zhcn.Range("G2").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/851f0658c8a9370b2bcf604134c416e9339c026a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf") | Out-Null

$zhcn.Range("F3").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md"
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/b3af6e214d11ac0ecaa6b80385723228932fa6be/e2e/430ce049-5402-4497-a1f8-90bcfc8e1e5c.md", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md") | Out-Null

$zhcn.Range("G3").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf"
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/851f0658c8a9370b2bcf604134c416e9339c026a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.zh-cn.xlf") | Out-Null

# Latest Handback DateTime (H) for zh-cn: "0001-01-01 00:00:00" -> "2016-03-12 14:49:08"
$zhcn.Range("H2").Value = "2016-03-12 14:49:08"
$zhcn.Range("H3").Value = "2016-03-12 14:49:08"

# de-de -------------------------------------------------------------------
$dede.Range("F2").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md"
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/b3af6e214d11ac0ecaa6b80385723228932fa6be/e2e/430ce049-5402-4497-a1f8-90bcfc8e1e5c.md", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md") | Out-Null

$dede.Range("G2").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43c6d7d750767f8ead252e5c8d7148a2ff97356f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf") | Out-Null

$dede.Range("F3").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md"
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/b3af6e214d11ac0ecaa6b80385723228932fa6be/e2e/430ce049-5402-4497-a1f8-90bcfc8e1e5c.md", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.md") | Out-Null

$dede.Range("G3").Value = "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf"
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43c6d7d750767f8ead252e5c8d7148a2ff97356f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "430ce049-5402-4497-a1f8-90bcfc8e1e5c.e430bfe384f31a6d9d916622e5e4e2a511fe8e0a.de-de.xlf") | Out-Null

# Latest Handback DateTime (H) for de-de: "0001-01-01 00:00:00" -> "2016-03-12 14:49:15"
$dede.Range("H2").Value = "2016-03-12 14:49:15"
$dede.Range("H3").Value = "2016-03-12 14:49:15"

"Handback report generated."
